# Apply the recorded edits to Barangay_Centers_Table.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix header typo: "longtitude" -> "longitude" (column D, row 1)
$ws.Range("D1").Value = "longitude"

# Update the "infected" counts in column A for the affected barangay rows
$ws.Range("A2").Value = 14
$ws.Range("A6").Value = 26
$ws.Range("A7").Value = 3
$ws.Range("A8").Value = 13
$ws.Range("A9").Value = 59
$ws.Range("A10").Value = 13
$ws.Range("A11").Value = 9
$ws.Range("A13").Value = 8
$ws.Range("A14").Value = 8
$ws.Range("A15").Value = 7
$ws.Range("A16").Value = 3
$ws.Range("A19").Value = 11
$ws.Range("A20").Value = 10
$ws.Range("A21").Value = 42
$ws.Range("A22").Value = 38
$ws.Range("A23").Value = 42
$ws.Range("A24").Value = 11
$ws.Range("A25").Value = 21
$ws.Range("A27").Value = 13
$ws.Range("A28").Value = 7
$ws.Range("A29").Value = 3
$ws.Range("A30").Value = 76
$ws.Range("A33").Value = 22
$ws.Range("A34").Value = 6
$ws.Range("A35").Value = 7
$ws.Range("A36").Value = 15
$ws.Range("A37").Value = 19
$ws.Range("A38").Value = 8
$ws.Range("A41").Value = 3
$ws.Range("A42").Value = 6
$ws.Range("A43").Value = 5

# Match the latest view state saved by Excel: zoomed in to 130%, G14 selected
$excel.ActiveWindow.Zoom = 130
$ws.Range("G14").Select()
